# "Add files via upload" edit for SKF_NU 234 ECM.xlsx
#
# On the active sheet (Лист1) a new "Subtype" attribute row is inserted right
# below the existing "Type" row, the "Type" value is shortened, and the
# leftover blank spacer row further down is removed so the sheet keeps the
# same overall number of rows. The active cell selection also moves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row right after the "Type" row (row 3) for "Subtype" ---
$ws.Rows("4:4").Insert()

# Reuse the formatting/styles of row 3 for the newly inserted row 4
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Match the row height used by the other data rows on the sheet
$ws.Rows("4:4").RowHeight = 15.75

# --- Populate the new "Subtype" row label ---
$ws.Range("A4").Value = "Subtype"

# --- Shorten the "Type" value in row 3 ---
$ws.Range("B3").Value = "Roller bearings"

# --- Populate the new "Subtype" row value ---
$ws.Range("B4").Value = "Single row cilindrical"

# --- Drop the now-redundant extra blank row so row count stays the same ---
$ws.Rows("17:17").Delete()

# --- Update the selected / active cell ---
$ws.Range("A17").Select()
